$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(6)
$tbl = $shape.Table
$cell = $tbl.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange

# Replace "prevAddressBook " (including its trailing space) with "prevCatalogue "
# so that the paragraph reads "prevCatalogue = s3" overall, matching the
# authors edit (renaming prevAddressBook -> prevCatalogue and moving the
# separating space onto the first run).
$sub = $tr.Characters(17, 16)
$sub.Text = "prevCatalogue "
